$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(4, 17).Value = "Phút tăng ca đêm"
Write-Host "set"
